$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = '41.927.44'
$ws.Range("E2").Value = '  -1.21%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = '2.210.89'
$ws.Range("E3").Value = '  -1.61%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = '  +0.19%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.17'
$ws.Range("E5").Value = '  -2.03%  '

# Row 6: 'XRP' -> 'XRP'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -1.05%  '

# Row 7: 'Solana' -> 'Solana'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.13'
$ws.Range("E7").Value = '  -4.65%  '

# Row 8: 'USDC' -> 'USDC'
$ws.Range("E8").Value = '  +0.14%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  -3.59%  '

# Row 10: 'Avalanche' -> 'Avalanche'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.57'
$ws.Range("E10").Value = '  -5.54%  '

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  -0.55%  '

# Row 12: 'Polkadot' -> 'TRON'
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  +0.30%  '

# Row 13: 'TRON' -> 'Polkadot'
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.94'
$ws.Range("E13").Value = '  -4.45%  '

# Row 14: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = '2.548.34'
$ws.Range("E14").Value = '  -1.43%  '

# Row 15: 'Chainlink' -> 'Chainlink'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.14'
$ws.Range("E15").Value = '  -2.89%  '

# Row 16: 'Polygon' -> 'Polygon'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.831'
$ws.Range("E16").Value = '  -2.89%  '

# Row 17: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D17").Value = '2.213.93'
$ws.Range("E17").Value = '  -1.62%  '

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D18").Value = '41.834.32'
$ws.Range("E18").Value = '  -1.06%  '

# Row 19: 'ShibaInu' -> 'ShibaInu'
$ws.Range("E19").Value = '  +2.17%  '

# Row 20: 'Uniswap' -> 'Litecoin'
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.33'
$ws.Range("E20").Value = '  +0.17%  '

# Row 21: 'Litecoin' -> 'Uniswap'
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  -0.72%  '

# Row 22: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.11'
$ws.Range("E22").Value = '  +20.02%  '

# Row 23: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.81'
$ws.Range("E23").Value = '  -1.27%  '

# Row 24: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  -8.98%  '

# Row 25: 'Dai' -> 'Dai'
$ws.Range("E25").Value = '  +0.19%  '

# Row 26: 'Cosmos' -> 'Cosmos'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.37'
$ws.Range("E26").Value = '  -1.04%  '

# Row 27: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range("E27").Value = '  +0.13%  '

# Row 28: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("E28").Value = '  -2.06%  '

# Row 29: 'Toncoin' -> 'Toncoin'
$ws.Range("E29").Value = '  -0.90%  '

# Row 30: 'Monero' -> 'Monero'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.07'
$ws.Range("E30").Value = '  -0.69%  '

# Row 31: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.35'
$ws.Range("E31").Value = '  -1.53%  '

# Row 32: 'Filecoin' -> 'Filecoin'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.57'
$ws.Range("E32").Value = '  +4.80%  '

# Row 33: 'Hedera' -> 'Hedera'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0792'
$ws.Range("E33").Value = '  -4.09%  '

# Row 34: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.77'
$ws.Range("E34").Value = '  -2.54%  '

# Row 35: 'Stellar' -> 'Stellar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.124'
$ws.Range("E35").Value = '  -0.60%  '

# Row 36: 'Kaspa' -> 'Kaspa'
$ws.Range("E36").Value = '  -12.29%  '

# Row 37: 'RenderToken' -> 'RenderToken'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.19'
$ws.Range("E37").Value = '  -7.69%  '

# Row 38: 'VeChain' -> 'VeChain'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0298'
$ws.Range("E38").Value = '  -5.70%  '

# Row 39: 'Celestia' -> 'Celestia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.48'
$ws.Range("E39").Value = '  -3.91%  '

# Row 40: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.10'
$ws.Range("E40").Value = '  -3.68%  '

# Row 41: 'MultiversX' -> 'THORChain'
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.61'
$ws.Range("E41").Value = '  -3.39%  '

# Row 42: 'THORChain' -> 'MultiversX'
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '63.82'
$ws.Range("E42").Value = '  -0.39%  '

# Row 43: 'Algorand' -> 'Algorand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.196'
$ws.Range("E43").Value = '  -3.24%  '

# Row 44: 'FraxShare' -> 'FraxShare'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.65'
$ws.Range("E44").Value = '  -1.66%  '

# Row 45: 'Aave' -> 'Aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.73'
$ws.Range("E45").Value = '  -4.98%  '

# Row 46: 'Cronos' -> 'Cronos'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0997'
$ws.Range("E46").Value = '  -2.43%  '

# Row 47: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.33'
$ws.Range("E47").Value = '  -0.38%  '

# Row 48: 'ARBITRUM' -> 'TrustWalletToken'
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.16'
$ws.Range("E48").Value = '  -1.96%  '

# Row 49: 'TrustWalletToken' -> 'ARBITRUM'
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.10'
$ws.Range("E49").Value = '  -2.67%  '

# Row 50: 'HuobiToken' -> 'HuobiToken'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.68'
$ws.Range("E50").Value = '  -0.85%  '

# Row 51: 'RocketPoolETH' -> 'RocketPoolETH'
$ws.Range("D51").Value = '2.423.62'
